$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Helper: write a literal text value (never auto-converted to a date/number)
# into a cell by routing it through a scratch cell's text-formula result and
# Copy/PasteSpecial, which preserves the literal string instead of Excel's
# "smart" type inference that .Value normally applies (e.g. "2025-11-16" ->
# a date serial). We reuse the same scratch cell for every distinct date.
function Set-LiteralText {
    param($cell, [string]$text)
    $scratch = $ws.Range("Z100")
    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy() | Out-Null
    $cell.PasteSpecial() | Out-Null
    $scratch.ClearContents() | Out-Null
}

# ---------------------------------------------------------------------------
# Row 7: existing row is repurposed from a "q" (quote) row into an "i"
# (invoice) row belonging to base_id 20251116-002.
# ---------------------------------------------------------------------------
$ws.Range("A7").Value = "20251116-002"
Set-LiteralText $ws.Range("B7") "2025-11-16"
$ws.Range("C7").Value = "i"
$ws.Range("D7").Value = "INV-20251116-003"
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = "gfgsfhwrth"
$ws.Range("G7").ClearContents() | Out-Null
$ws.Range("H7").Value = "Abu Dhabi - Al Shamkha"
$ws.Range("I7").Value = "QUO-20251116-004"

# ---------------------------------------------------------------------------
# New rows 8-17
# ---------------------------------------------------------------------------

# Row 8
$ws.Range("A8").Value = "20251116-002"
Set-LiteralText $ws.Range("B8") "2025-11-16"
$ws.Range("C8").Value = "i"
$ws.Range("D8").Value = "INV-20251116-004"
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = "Fahad Ahmed Mohammed"
$ws.Range("G8").ClearContents() | Out-Null
$ws.Range("H8").Value = "Abu Dhabi - Al Shamkha"
$ws.Range("I8").Value = "QUO-20251116-004"

# Row 9
$ws.Range("A9").Value = "20251116-002"
Set-LiteralText $ws.Range("B9") "2025-11-16"
$ws.Range("C9").Value = "i"
$ws.Range("D9").Value = "INV-20251116-005"
$ws.Range("E9").Value = 1100
$ws.Range("F9").ClearContents() | Out-Null
$ws.Range("G9").ClearContents() | Out-Null
$ws.Range("H9").Value = "Abu Dhabi - Al Shamkha"
$ws.Range("I9").Value = "QUO-20251116-004"

# Row 10
$ws.Range("A10").Value = "20251116-002"
Set-LiteralText $ws.Range("B10") "2025-11-16"
$ws.Range("C10").Value = "i"
$ws.Range("D10").Value = "INV-20251116-006"
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = "Fahad Ahmed Mohammed"
$ws.Range("G10").ClearContents() | Out-Null
$ws.Range("H10").Value = "Abu Dhabi - Al Shamkha"
$ws.Range("I10").Value = "QUO-20251116-004"

# Row 11
$ws.Range("A11").Value = "20251116-002"
Set-LiteralText $ws.Range("B11") "2025-11-16"
$ws.Range("C11").Value = "i"
$ws.Range("D11").Value = "INV-20251116-007"
$ws.Range("E11").Value = 0
$ws.Range("F11").ClearContents() | Out-Null
$ws.Range("G11").ClearContents() | Out-Null
$ws.Range("H11").Value = "Abu Dhabi - Al Shamkha"
$ws.Range("I11").Value = "QUO-20251116-004"

# Row 12
$ws.Range("A12").Value = "20251116-012"
Set-LiteralText $ws.Range("B12") "2025-11-16"
$ws.Range("C12").Value = "q"
$ws.Range("D12").Value = "QUO-20251116-001"
$ws.Range("E12").Value = 0
$ws.Range("F12").ClearContents() | Out-Null
$ws.Range("G12").ClearContents() | Out-Null
$ws.Range("H12").Value = "Abu Dhabi - Al Shamkha"
$ws.Range("I12").ClearContents() | Out-Null

# Row 13
$ws.Range("A13").Value = "20251116-002"
Set-LiteralText $ws.Range("B13") "2025-11-17"
$ws.Range("C13").Value = "i"
$ws.Range("D13").Value = "INV-20251117-008"
$ws.Range("E13").Value = 2200
$ws.Range("F13").Value = "Fahad Ahmed Mohammed"
$ws.Range("G13").ClearContents() | Out-Null
$ws.Range("H13").Value = "Abu Dhabi - Al Shamkha"
$ws.Range("I13").Value = "QUO-20251116-004"

# Row 14
$ws.Range("A14").Value = "20251117-003"
Set-LiteralText $ws.Range("B14") "2025-11-17"
$ws.Range("C14").Value = "q"
$ws.Range("D14").Value = "QUO-20251117-004"
$ws.Range("E14").Value = 3300
$ws.Range("F14").ClearContents() | Out-Null
$ws.Range("G14").ClearContents() | Out-Null
$ws.Range("H14").Value = "Abu Dhabi - Al Shamkha"
$ws.Range("I14").Value = "PDF"

# Row 15
$ws.Range("A15").Value = "20251117-003"
Set-LiteralText $ws.Range("B15") "2025-11-17"
$ws.Range("C15").Value = "q"
$ws.Range("D15").Value = "QUO-20251117-001"
$ws.Range("E15").Value = 0
$ws.Range("F15").ClearContents() | Out-Null
$ws.Range("G15").ClearContents() | Out-Null
$ws.Range("H15").Value = "Abu Dhabi - Al Shamkha"
$ws.Range("I15").ClearContents() | Out-Null

# Row 16
$ws.Range("A16").Value = "20251116-002"
Set-LiteralText $ws.Range("B16") "2025-11-18"
$ws.Range("C16").Value = "i"
$ws.Range("D16").Value = "INV-20251118-009"
$ws.Range("E16").Value = 949.93
$ws.Range("F16").Value = "Fahad Ahmed Mohammed"
$ws.Range("G16").ClearContents() | Out-Null
$ws.Range("H16").Value = "Abu Dhabi - Al Shamkha"
$ws.Range("I16").Value = "QUO-20251116-004"

# Row 17
$ws.Range("A17").Value = "20251118-002"
Set-LiteralText $ws.Range("B17") "2025-11-18"
$ws.Range("C17").Value = "q"
$ws.Range("D17").Value = "QUO-20251118-003"
$ws.Range("E17").Value = 1230
$ws.Range("F17").ClearContents() | Out-Null
$ws.Range("G17").ClearContents() | Out-Null
$ws.Range("H17").Value = "Abu Dhabi - Al Shamkha"
$ws.Range("I17").Value = "PDF"
